$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format before writing, so purely-numeric-
# looking strings like "1.003" or "321.80" are stored as text (matching the
# original inlineStr cells) instead of being auto-parsed into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.430.28'
$ws.Range('E2').Value = '  -2.91%  '
$ws.Range('D3').Value = '1.743.81'
$ws.Range('E3').Value = '  -3.43%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '321.80'
$ws.Range('E5').Value = '  -4.30%  '
$ws.Range('D7').Value = '0.4233'
$ws.Range('E7').Value = '  -8.37%  '
$ws.Range('D8').Value = '0.3584'
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = '45.47'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').Value = '0.07431'
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('E11').Value = '  -3.55%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '21.45'
$ws.Range('E13').Value = '  -4.40%  '
$ws.Range('D14').Value = '6.117'
$ws.Range('E14').Value = '  -3.60%  '
$ws.Range('D15').Value = '7.192'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '1.743.04'
$ws.Range('E16').Value = '  -3.47%  '
$ws.Range('D17').Value = '0.00001066'
$ws.Range('E17').Value = '  -3.04%  '
$ws.Range('D18').Value = '87.97'
$ws.Range('E18').Value = '  +7.37%  '
$ws.Range('D19').Value = '0.06113'
$ws.Range('E19').Value = '  -9.02%  '
$ws.Range('D21').Value = '16.87'
$ws.Range('E21').Value = '  -3.66%  '
$ws.Range('D22').Value = '6.104'
$ws.Range('E22').Value = '  -5.01%  '
$ws.Range('D23').Value = '0.5245'
$ws.Range('E23').Value = '  -6.18%  '
$ws.Range('D24').Value = '27.464.38'
$ws.Range('E24').Value = '  -2.77%  '
$ws.Range('E25').Value = '  -3.47%  '
$ws.Range('D26').Value = '2.340'
$ws.Range('E26').Value = '  -2.86%  '
$ws.Range('D27').Value = '20.39'
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('D28').Value = '2.380'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').Value = '152.55'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('D30').Value = '1.938.81'
$ws.Range('E30').Value = '  -3.61%  '
$ws.Range('D31').Value = '126.05'
$ws.Range('E31').Value = '  -5.59%  '
$ws.Range('D32').Value = '1.199'
$ws.Range('E32').Value = '  -4.79%  '
$ws.Range('D33').Value = '5.674'
$ws.Range('E33').Value = '  -3.38%  '
$ws.Range('D34').Value = '0.09113'
$ws.Range('E34').Value = '  -4.83%  '
$ws.Range('D35').Value = '3.624'
$ws.Range('E35').Value = '  -10.18%  '
$ws.Range('D36').Value = '12.66'
$ws.Range('E36').Value = '  +4.33%  '
$ws.Range('D37').Value = '0.02296'
$ws.Range('E37').Value = '  -2.63%  '
$ws.Range('D38').Value = '0.2138'
$ws.Range('E38').Value = '  -3.90%  '
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').Value = '0.6394'
$ws.Range('E41').Value = '  -3.95%  '
$ws.Range('D42').Value = '1.188'
$ws.Range('E42').Value = '  -3.96%  '
$ws.Range('D43').Value = '1.424'
$ws.Range('E43').Value = '  -4.85%  '
$ws.Range('D45').Value = '7.897'
$ws.Range('E45').Value = '  -3.44%  '
$ws.Range('D46').Value = '13.68'
$ws.Range('E46').Value = '  -4.16%  '
$ws.Range('D47').Value = '3.712'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('D48').Value = '0.5874'
$ws.Range('E48').Value = '  -4.38%  '
$ws.Range('D49').Value = '125.15'
$ws.Range('E49').Value = '  -3.74%  '
$ws.Range('D50').Value = '1.945'
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('D51').Value = '0.06833'
$ws.Range('E51').Value = '  -4.57%  '

# Restore the default (Normal) style so the cells end up with no explicit
# style index, matching the original workbook formatting.
$priceRange.Style = "Normal"
